$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.803.29'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.52%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.111.79'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.11%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.88'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.00'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.79%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.108.87'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.521'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.31%  '
$ws.Range("E10").Value = '  -2.37%  '
$ws.Range("E11").Value = '  -0.63%  '
$ws.Range("E12").Value = '  +0.29%  '
$ws.Range("E13").Value = '  -1.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.01'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.79%  '
$ws.Range("E15").Value = '  -1.20%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.628.33'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.12%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.783.51'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.29%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.18'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.51%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.110.87'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.32'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.17%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '476.00'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.95%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.714'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.88'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.80%  '
$ws.Range("E24").Value = '  +3.90%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.10'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.27%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.10'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.01%  '
$ws.Range("E28").Value = '  +0.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.87'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.19%  '
$ws.Range("E30").Value = '  -0.78%  '
$ws.Range("E31").Value = '  +0.61%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.56'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.46%  '
$ws.Range("E33").Value = '  +1.60%  '
$ws.Range("E34").Value = '  -7.35%  '
$ws.Range("E35").Value = '  -0.13%  '
$ws.Range("E36").Value = '  -0.19%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.973'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.88%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '46.88'
$ws.Range("D38").Style = "Normal"
$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '50.14'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.08%  '
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.06'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.47%  '
$ws.Range("E41").Value = '  -2.23%  '
$ws.Range("E42").Value = '  +0.05%  '
$ws.Range("E43").Value = '  +0.50%  '
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.830.70'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.61%  '
$ws.Range("B45").Value = 'Bittensor'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '386.41'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.17%  '
$ws.Range("B46").Value = 'dogwifhat'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.59'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -7.89%  '
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0355'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.42%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '135.55'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.55%  '
$ws.Range("E49").Value = '  +0.00%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '24.86'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.14%  '
$ws.Range("E51").Value = '  -1.67%  '
